# Update column F (dSF) values for a handful of rows, per the "repull data,
# push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F6").Value = -8
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = 3
$ws.Range("F12").Value = 2
$ws.Range("F15").Value = -1
$ws.Range("F18").Value = -6
$ws.Range("F23").Value = -3
